$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $st = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $st
}

$ws.Range('D2').Value = '68.493.50'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '2.454.61'
$ws.Range('E3').Value = '  -2.13%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue 'D5' '562.93'
$ws.Range('E5').Value = '  -2.26%  '
Set-TextValue 'D6' '163.22'
$ws.Range('E6').Value = '  -2.30%  '
$ws.Range('E7').Value = '  +0.00%  '
Set-TextValue 'D8' '0.504'
$ws.Range('E8').Value = '  -2.21%  '
$ws.Range('D9').Value = '2.453.01'
$ws.Range('E9').Value = '  -2.18%  '
Set-TextValue 'D10' '0.151'
$ws.Range('E10').Value = '  -6.07%  '
$ws.Range('E11').Value = '  -1.92%  '
Set-TextValue 'D12' '0.342'
$ws.Range('E12').Value = '  -5.30%  '
Set-TextValue 'D13' '4.80'
$ws.Range('E13').Value = '  -2.79%  '
$ws.Range('D14').Value = '2.905.22'
$ws.Range('E14').Value = '  -2.12%  '
$ws.Range('D15').Value = '68.296.29'
$ws.Range('E15').Value = '  -1.67%  '
Set-TextValue 'D16' '0.0000171'
$ws.Range('E16').Value = '  -3.86%  '
Set-TextValue 'D17' '23.73'
$ws.Range('E17').Value = '  -4.56%  '
$ws.Range('D18').Value = '2.449.07'
$ws.Range('E18').Value = '  -2.33%  '
Set-TextValue 'D19' '11.00'
$ws.Range('E19').Value = '  -2.71%  '
Set-TextValue 'D20' '350.29'
$ws.Range('E20').Value = '  -0.38%  '
Set-TextValue 'D21' '7.21'
$ws.Range('E21').Value = '  -5.00%  '
Set-TextValue 'D22' '3.81'
$ws.Range('E22').Value = '  -2.88%  '
$ws.Range('E23').Value = '  +0.03%  '
Set-TextValue 'D24' '1.86'
$ws.Range('E24').Value = '  -5.33%  '
Set-TextValue 'D25' '68.06'
$ws.Range('E25').Value = '  -3.66%  '
Set-TextValue 'D26' '3.76'
$ws.Range('E26').Value = '  -4.86%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D27' '1.05'
$ws.Range('E27').Value = '  +4.39%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.574.07'
$ws.Range('E28').Value = '  -3.45%  '
Set-TextValue 'D29' '8.28'
$ws.Range('E29').Value = '  -6.57%  '
$ws.Range('D30').Value = '0.0₃0842'
$ws.Range('E30').Value = '  -5.68%  '
Set-TextValue 'D31' '7.34'
$ws.Range('E31').Value = '  -6.77%  '
Set-TextValue 'D32' '0.999'
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('E33').Value = '  -4.36%  '
Set-TextValue 'D34' '430.51'
$ws.Range('E34').Value = '  -6.65%  '
Set-TextValue 'D35' '1.68'
$ws.Range('E35').Value = '  -3.37%  '
Set-TextValue 'D36' '3.03'
$ws.Range('E36').Value = '  +104.08%  '
Set-TextValue 'D37' '157.06'
$ws.Range('E37').Value = '  -1.77%  '
Set-TextValue 'D38' '19.01'
$ws.Range('E38').Value = '  -0.33%  '
$ws.Range('B39').Value = 'USDe'
$ws.Range('C39').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D39' '1.00'
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D40' '0.110'
$ws.Range('E40').Value = '  -5.60%  '
Set-TextValue 'D41' '17.94'
$ws.Range('E41').Value = '  -2.95%  '
Set-TextValue 'D42' '0.306'
$ws.Range('E42').Value = '  -4.09%  '
Set-TextValue 'D43' '4.51'
$ws.Range('E43').Value = '  -3.87%  '
$ws.Range('E44').Value = '  -4.34%  '
Set-TextValue 'D45' '1.09'
$ws.Range('E45').Value = '  -0.81%  '
Set-TextValue 'D46' '2.08'
$ws.Range('E46').Value = '  -6.59%  '
Set-TextValue 'D47' '135.97'
$ws.Range('E47').Value = '  -4.47%  '
Set-TextValue 'D48' '3.36'
$ws.Range('E48').Value = '  -3.42%  '
Set-TextValue 'D49' '0.490'
$ws.Range('E49').Value = '  -5.91%  '
Set-TextValue 'D50' '0.0715'
$ws.Range('E50').Value = '  -2.67%  '
Set-TextValue 'D51' '0.563'
$ws.Range('E51').Value = '  -2.51%  '
